$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 6.603177
$ws.Range("H2").Value2 = 19.809531
$ws.Range("I2").Value2 = 0.5135477412645301
$ws.Range("J2").Value2 = 0.5135477412645302
$ws.Range("M2").Value2 = 0.034325
$ws.Range("N2").Value2 = 0.102975
$ws.Range("O2").Value2 = 0.004508979075184418
$ws.Range("P2").Value2 = 0.004508979075184418
$ws.Range("Q2").Value2 = 0.226654050525
$ws.Range("R2").Value2 = 2.039886454725
$ws.Range("S2").Value2 = 0.002315576019469988
$ws.Range("T2").Value2 = 0.002315576019469988

$ws.Range("G3").Value2 = 6.603177
$ws.Range("H3").Value2 = 19.809531
$ws.Range("I3").Value2 = 0.5135477412645301
$ws.Range("J3").Value2 = 0.5135477412645302
$ws.Range("O3").Value2 = 0.9142039036746329
$ws.Range("P3").Value2 = 0.9142039036746329
$ws.Range("Q3").Value2 = 45.954530797007
$ws.Range("R3").Value2 = 413.590777173063
$ws.Range("S3").Value2 = 0.4694873497873238
$ws.Range("T3").Value2 = 0.4694873497873238

$ws.Range("G4").Value2 = 6.603177
$ws.Range("H4").Value2 = 19.809531
$ws.Range("I4").Value2 = 0.5135477412645301
$ws.Range("J4").Value2 = 0.5135477412645302
$ws.Range("M4").Value2 = 0.5818573333333333
$ws.Range("N4").Value2 = 1.745572
$ws.Range("O4").Value2 = 0.0764335772976724
$ws.Range("P4").Value2 = 0.0764335772976724
$ws.Range("Q4").Value2 = 3.842106960748
$ws.Range("R4").Value2 = 34.578962646732
$ws.Range("S4").Value2 = 0.03925229097798753
$ws.Range("T4").Value2 = 0.03925229097798753

$ws.Range("G5").Value2 = 6.603177
$ws.Range("H5").Value2 = 19.809531
$ws.Range("I5").Value2 = 0.5135477412645301
$ws.Range("J5").Value2 = 0.5135477412645302
$ws.Range("M5").Value2 = 0.036948
$ws.Range("N5").Value2 = 0.110844
$ws.Range("O5").Value2 = 0.004853539952510238
$ws.Range("P5").Value2 = 0.004853539952510237
$ws.Range("Q5").Value2 = 0.243974183796
$ws.Range("R5").Value2 = 2.195767654164
$ws.Range("S5").Value2 = 0.002492524479748787
$ws.Range("T5").Value2 = 0.002492524479748787

$ws.Range("I6").Value2 = 0.02944398858046029
$ws.Range("J6").Value2 = 0.0294439885804603
$ws.Range("M6").Value2 = 0.034325
$ws.Range("N6").Value2 = 0.102975
$ws.Range("O6").Value2 = 0.004508979075184418
$ws.Range("P6").Value2 = 0.004508979075184418
$ws.Range("Q6").Value2 = 0.01299509030833333
$ws.Range("R6").Value2 = 0.116955812775
$ws.Range("S6").Value2 = 0.0001327623283992644
$ws.Range("T6").Value2 = 0.0001327623283992644

$ws.Range("I7").Value2 = 0.02944398858046029
$ws.Range("J7").Value2 = 0.0294439885804603
$ws.Range("O7").Value2 = 0.9142039036746329
$ws.Range("P7").Value2 = 0.9142039036746329
$ws.Range("S7").Value2 = 0.02691780930000811
$ws.Range("T7").Value2 = 0.02691780930000811

$ws.Range("I8").Value2 = 0.02944398858046029
$ws.Range("J8").Value2 = 0.0294439885804603
$ws.Range("M8").Value2 = 0.5818573333333333
$ws.Range("N8").Value2 = 1.745572
$ws.Range("O8").Value2 = 0.0764335772976724
$ws.Range("P8").Value2 = 0.0764335772976724
$ws.Range("Q8").Value2 = 0.2202851738742222
$ws.Range("R8").Value2 = 1.982566564868
$ws.Range("S8").Value2 = 0.002250509377116395
$ws.Range("T8").Value2 = 0.002250509377116395

$ws.Range("I9").Value2 = 0.02944398858046029
$ws.Range("J9").Value2 = 0.0294439885804603
$ws.Range("M9").Value2 = 0.036948
$ws.Range("N9").Value2 = 0.110844
$ws.Range("O9").Value2 = 0.004853539952510238
$ws.Range("P9").Value2 = 0.004853539952510237
$ws.Range("Q9").Value2 = 0.013988131004
$ws.Range("R9").Value2 = 0.125893179036
$ws.Range("S9").Value2 = 0.0001429075749365192
$ws.Range("T9").Value2 = 0.0001429075749365192

$ws.Range("G10").Value2 = 3.441487333333333
$ws.Range("H10").Value2 = 10.324462
$ws.Range("I10").Value2 = 0.2676541983690312
$ws.Range("J10").Value2 = 0.2676541983690313
$ws.Range("M10").Value2 = 0.034325
$ws.Range("N10").Value2 = 0.102975
$ws.Range("O10").Value2 = 0.004508979075184418
$ws.Range("P10").Value2 = 0.004508979075184418
$ws.Range("Q10").Value2 = 0.1181290527166667
$ws.Range("R10").Value2 = 1.06316147445
$ws.Range("S10").Value2 = 0.001206847179831221
$ws.Range("T10").Value2 = 0.001206847179831221

$ws.Range("G11").Value2 = 3.441487333333333
$ws.Range("H11").Value2 = 10.324462
$ws.Range("I11").Value2 = 0.2676541983690312
$ws.Range("J11").Value2 = 0.2676541983690313
$ws.Range("O11").Value2 = 0.9142039036746329
$ws.Range("P11").Value2 = 0.9142039036746329
$ws.Range("Q11").Value2 = 23.95088540670289
$ws.Range("R11").Value2 = 215.557968660326
$ws.Range("S11").Value2 = 0.2446905129838729
$ws.Range("T11").Value2 = 0.2446905129838729

$ws.Range("G12").Value2 = 3.441487333333333
$ws.Range("H12").Value2 = 10.324462
$ws.Range("I12").Value2 = 0.2676541983690312
$ws.Range("J12").Value2 = 0.2676541983690313
$ws.Range("M12").Value2 = 0.5818573333333333
$ws.Range("N12").Value2 = 1.745572
$ws.Range("O12").Value2 = 0.0764335772976724
$ws.Range("P12").Value2 = 0.0764335772976724
$ws.Range("Q12").Value2 = 2.002454642473778
$ws.Range("R12").Value2 = 18.022091782264
$ws.Range("S12").Value2 = 0.02045776786008589
$ws.Range("T12").Value2 = 0.02045776786008589

$ws.Range("G13").Value2 = 3.441487333333333
$ws.Range("H13").Value2 = 10.324462
$ws.Range("I13").Value2 = 0.2676541983690312
$ws.Range("J13").Value2 = 0.2676541983690313
$ws.Range("M13").Value2 = 0.036948
$ws.Range("N13").Value2 = 0.110844
$ws.Range("O13").Value2 = 0.004853539952510238
$ws.Range("P13").Value2 = 0.004853539952510237
$ws.Range("Q13").Value2 = 0.127156073992
$ws.Range("R13").Value2 = 1.144404665928
$ws.Range("S13").Value2 = 0.001299070345241194
$ws.Range("T13").Value2 = 0.001299070345241194

$ws.Range("G14").Value2 = 2.434707333333333
$ws.Range("H14").Value2 = 7.304122
$ws.Range("I14").Value2 = 0.1893540717859783
$ws.Range("J14").Value2 = 0.1893540717859783
$ws.Range("M14").Value2 = 0.034325
$ws.Range("N14").Value2 = 0.102975
$ws.Range("O14").Value2 = 0.004508979075184418
$ws.Range("P14").Value2 = 0.004508979075184418
$ws.Range("Q14").Value2 = 0.08357132921666667
$ws.Range("R14").Value2 = 0.75214196295
$ws.Range("S14").Value2 = 0.0008537935474839444
$ws.Range("T14").Value2 = 0.0008537935474839445

$ws.Range("G15").Value2 = 2.434707333333333
$ws.Range("H15").Value2 = 7.304122
$ws.Range("I15").Value2 = 0.1893540717859783
$ws.Range("J15").Value2 = 0.1893540717859783
$ws.Range("O15").Value2 = 0.9142039036746329
$ws.Range("P15").Value2 = 0.9142039036746329
$ws.Range("Q15").Value2 = 16.94424261705622
$ws.Range("R15").Value2 = 152.498183553506
$ws.Range("S15").Value2 = 0.173108231603428
$ws.Range("T15").Value2 = 0.1731082316034281

$ws.Range("G16").Value2 = 2.434707333333333
$ws.Range("H16").Value2 = 7.304122
$ws.Range("I16").Value2 = 0.1893540717859783
$ws.Range("J16").Value2 = 0.1893540717859783
$ws.Range("M16").Value2 = 0.5818573333333333
$ws.Range("N16").Value2 = 1.745572
$ws.Range("O16").Value2 = 0.0764335772976724
$ws.Range("P16").Value2 = 0.0764335772976724
$ws.Range("Q16").Value2 = 1.416652316420445
$ws.Range("R16").Value2 = 12.749870847784
$ws.Range("S16").Value2 = 0.01447300908248258
$ws.Range("T16").Value2 = 0.01447300908248258

$ws.Range("G17").Value2 = 2.434707333333333
$ws.Range("H17").Value2 = 7.304122
$ws.Range("I17").Value2 = 0.1893540717859783
$ws.Range("J17").Value2 = 0.1893540717859783
$ws.Range("M17").Value2 = 0.036948
$ws.Range("N17").Value2 = 0.110844
$ws.Range("O17").Value2 = 0.004853539952510238
$ws.Range("P17").Value2 = 0.004853539952510237
$ws.Range("Q17").Value2 = 0.089957566552
$ws.Range("R17").Value2 = 0.8096180989680001
$ws.Range("S17").Value2 = 0.0009190375525837373
$ws.Range("T17").Value2 = 0.0009190375525837373
